$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to remain Text,
# matching the source workbook (all data cells are stored as text),
# then strip the temporary Text number-format back to the default style
# so no stray formatting is left behind.
function Set-TextCell($ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "43.938.01"
Set-TextCell "E2" "  -0.08%  "

Set-TextCell "D3" "2.234.06"
Set-TextCell "E3" "  -1.00%  "

Set-TextCell "D4" "1.00"
Set-TextCell "E4" "  +0.16%  "

Set-TextCell "D5" "303.87"
Set-TextCell "E5" "  -4.79%  "

Set-TextCell "D6" "94.58"
Set-TextCell "E6" "  -7.24%  "

Set-TextCell "D7" "0.567"
Set-TextCell "E7" "  -1.73%  "

Set-TextCell "E8" "  +0.18%  "

Set-TextCell "D9" "0.518"
Set-TextCell "E9" "  -6.57%  "

Set-TextCell "D10" "34.41"
Set-TextCell "E10" "  -7.97%  "

Set-TextCell "D11" "0.0804"
Set-TextCell "E11" "  -3.40%  "

Set-TextCell "D12" "7.14"
Set-TextCell "E12" "  -6.54%  "

Set-TextCell "E13" "  -2.89%  "

Set-TextCell "D14" "2.574.29"
Set-TextCell "E14" "  -1.02%  "

Set-TextCell "D15" "2.318.60"
Set-TextCell "E15" "  +3.00%  "

Set-TextCell "D16" "0.811"
Set-TextCell "E16" "  -5.80%  "

Set-TextCell "D17" "13.36"
Set-TextCell "E17" "  -7.82%  "

Set-TextCell "D18" "43.775.95"
Set-TextCell "E18" "  -0.25%  "

Set-TextCell "D19" "0.0₃0951"
Set-TextCell "E19" "  -3.74%  "

Set-TextCell "D20" "12.13"
Set-TextCell "E20" "  -9.97%  "

Set-TextCell "D21" "6.13"
Set-TextCell "E21" "  -6.32%  "

Set-TextCell "D22" "64.43"
Set-TextCell "E22" "  -2.30%  "

Set-TextCell "D23" "236.64"
Set-TextCell "E23" "  +0.42%  "

Set-TextCell "D24" "2.90"
Set-TextCell "E24" "  -7.51%  "

Set-TextCell "D25" "1.01"
Set-TextCell "E25" "  +0.51%  "

Set-TextCell "D26" "1.93"
Set-TextCell "E26" "  -8.49%  "

Set-TextCell "D27" "9.80"
Set-TextCell "E27" "  -4.31%  "

Set-TextCell "E28" "  -2.68%  "

Set-TextCell "D29" "36.04"
Set-TextCell "E29" "  -3.76%  "

Set-TextCell "D30" "19.97"
Set-TextCell "E30" "  -1.15%  "

Set-TextCell "D31" "5.84"
Set-TextCell "E31" "  -6.33%  "

Set-TextCell "D32" "152.06"
Set-TextCell "E32" "  -4.98%  "

Set-TextCell "D33" "0.0801"
Set-TextCell "E33" "  -6.18%  "

Set-TextCell "B34" "LidoDAOToken"
Set-TextCell "C34" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D34" "3.24"
Set-TextCell "E34" "  +5.38%  "

Set-TextCell "B35" "WEMIXToken"
Set-TextCell "C35" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D35" "2.62"
Set-TextCell "E35" "  -2.66%  "

Set-TextCell "E36" "  -5.83%  "

Set-TextCell "D37" "0.118"
Set-TextCell "E37" "  -1.40%  "

Set-TextCell "E38" "  -9.70%  "

Set-TextCell "D39" "14.69"
Set-TextCell "E39" "  -10.38%  "

Set-TextCell "D40" "3.30"
Set-TextCell "E40" "  -11.62%  "

Set-TextCell "D41" "3.75"
Set-TextCell "E41" "  -11.06%  "

Set-TextCell "D42" "0.0296"
Set-TextCell "E42" "  -6.02%  "

Set-TextCell "E43" "  +0.16%  "

Set-TextCell "D44" "1.724.33"
Set-TextCell "E44" "  -4.36%  "

Set-TextCell "D45" "84.24"
Set-TextCell "E45" "  +1.89%  "

Set-TextCell "B46" "Algorand"
Set-TextCell "C46" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell "D46" "0.184"
Set-TextCell "E46" "  -7.24%  "

Set-TextCell "B47" "Aave"
Set-TextCell "C47" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D47" "99.16"
Set-TextCell "E47" "  -5.29%  "

Set-TextCell "D48" "4.84"
Set-TextCell "E48" "  -7.36%  "

Set-TextCell "D49" "14.44"
Set-TextCell "E49" "  +2.13%  "

Set-TextCell "D50" "7.98"
Set-TextCell "E50" "  -4.32%  "

Set-TextCell "D51" "68.11"
Set-TextCell "E51" "  -10.29%  "
